$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E:E").Delete()
